# Refresh the crypto price/volume snapshot (columns D and E) to the latest
# scrape, as produced by the scheduled GitHub Actions update job.
# D = Price text, E = Volume(1h) text (padded "  +x.xx%  " strings).
# Some price strings look numeric (e.g. "7.99", "0.0780"); force those
# particular cells to Text format first so Excel doesn't silently coerce
# them to floating point numbers and drop significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.016.06"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "3.410.96"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.72"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.18"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.99"
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("D12").Value = "3.996.69"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.36"
$ws.Range("E14").Value = "  -4.73%  "
$ws.Range("D15").Value = "3.396.76"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "62.000.09"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.39"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.48"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.96"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.55"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.80"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "3.573.72"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.06"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("E36").Value = "  +3.91%  "
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "169.36"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "30.36"
$ws.Range("E39").Value = "  -4.79%  "
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0780"
$ws.Range("E41").Value = "  +3.72%  "
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.35"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("D47").Value = "2.541.61"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("E48").Value = "  +2.63%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("E51").Value = "  +0.01%  "
